# feat(logic), #27: Add OR, XOR, AND, NOT functions.
#
# Mirrors the worked "circular.xlsx" example used by the formulas test
# suite: row 12 gains a new OR()-based circular chain (A12..E12) and a
# fresh blank row 13 is appended below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New example row demonstrating OR() inside the same circular-reference
# pattern already used by the rows above (B->C->D->B, E = B+1).
$ws.Range("A12").Value = $true
$ws.Range("B12").Formula = "=OR(A12,C12)"
$ws.Range("C12").Formula = "=D12"
$ws.Range("D12").Formula = "=B12"
$ws.Range("E12").Formula = "=B12+1"

# Give the new row the same centered style used by the header cells
# (B1:D1) instead of the plain "E12 only" style that used to sit there.
$ws.Range("A12:E12").HorizontalAlignment = -4108

# A fresh, blank row underneath it -- also centered -- so the sheet's
# used range grows to A1:E13.
$ws.Range("A13:E13").HorizontalAlignment = -4108

$ws.Range("B12").Select() | Out-Null

Write-Output "applied OR/XOR/AND/NOT sample row"
